$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timesheet")

# --- Update existing task durations -----------------------------------
$ws.Range("C4").Value2  = 3
$ws.Range("C5").Value2  = 1.5
$ws.Range("C7").Value2  = 3.25
$ws.Range("C10").Value2 = 5.5
$ws.Range("C11").Value2 = 3.75

# --- Remove the blank spacer row (old row 12) --------------------------
# This shifts the "Milestone 1/2/3" rows up by one (old 13/14/15 -> 12/13/14).
$ws.Rows("12:12").Delete()

# --- Insert 3 new task rows right after the (now shifted) Milestone 1 row
$ws.Rows("13:15").Insert()

# Copy the formatting (date/task/time column styles) from an existing
# task row (row 3) onto the freshly inserted rows so they pick up the
# same cell styles used throughout the table.
$ws.Range("A3:C3").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Range("A3:C3").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A3:C3").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fill in the new rows -----------------------------------------------
# Shared-string table order follows entry order: "Analyse + dossier" was
# typed first, then "Aanpassingen dossier", then "Analyse" - so we set
# the text values in that order to reproduce the same sharedStrings.xml
# layout, and fill in the rest (dates / hours) afterwards.
$ws.Range("B15").Value = "Analyse + dossier"
$ws.Range("B13").Value = "Aanpassingen dossier"
$ws.Range("B14").Value = "Analyse"

$ws.Range("A13").Value2 = 42390
$ws.Range("C13").Value2 = 1.5

$ws.Range("A14").Value2 = 42665
$ws.Range("C14").Value2 = 3.75

$ws.Range("A15").Value2 = 42666
$ws.Range("C15").Value2 = 1.75

# --- Selection / view state ---------------------------------------------
$ws.Range("C15").Select()
$excel.ActiveWindow.WindowState = $excel.ActiveWindow.WindowState
$wb.Windows.Item(1).Width = 17600
$wb.Windows.Item(1).Height = 16140
